$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.363.60'
$ws.Range('E2').Value = '  -1.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.564.31'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.15'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.500'
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '21.89'
$ws.Range('E8').Value = '  -1.96%  '
$ws.Range('E9').Value = '  -2.01%  '
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0867'
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.788.70'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.562.54'
$ws.Range('E13').Value = '  -1.31%  '
$ws.Range('E14').Value = '  -1.07%  '
$ws.Range('E15').Value = '  -3.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.36'
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.380.40'
$ws.Range('E17').Value = '  -0.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0688'
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '211.84'
$ws.Range('E19').Value = '  -2.69%  '
$ws.Range('E20').Value = '  -0.97%  '
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('E22').Value = '  -0.84%  '
$ws.Range('E23').Value = '  -0.33%  '
$ws.Range('E24').Value = '  +1.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.28'
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.71'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '14.97'
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('E30').Value = '  -0.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0469'
$ws.Range('E31').Value = '  +1.09%  '
$ws.Range('E32').Value = '  -0.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.363.92'
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('E35').Value = '  +1.37%  '
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('E38').Value = '  +0.63%  '
$ws.Range('E39').Value = '  -1.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.821'
$ws.Range('E40').Value = '  +0.59%  '
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('E42').Value = '  -0.45%  '
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.95'
$ws.Range('E44').Value = '  +0.56%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.28'
$ws.Range('E45').Value = '  +1.33%  '
$ws.Range('E46').Value = '  -1.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.700.56'
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.47'
$ws.Range('E48').Value = '  -2.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₇0993'
$ws.Range('E49').Value = '  -0.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0954'
$ws.Range('E50').Value = '  -2.16%  '
$ws.Range('E51').Value = '  -0.79%  '
